$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.945.20"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.263.66"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  -0.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "2.606.26"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "2.264.20"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.804"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "44.828.86"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +16.53%  "
$ws.Range("D20").Value = "0.0₃0916"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0788"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.38%  "
$ws.Range("D45").Value = "1.779.50"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "69.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.51%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  -0.61%  "
